$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 842, shifting rows 842:916
# down to 844:918 (weekly update: a new week of observations is added).
$ws.Rows("842:843").Insert()

# New row 842: "Primera" quality, "$/caja 36 atados" unit, date 2022-07-27
$ws.Range("A842").Value = 6
$ws.Range("B842").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C842").Value = "Metropolitana"
$ws.Range("D842").Value = "2022-07-27"
$ws.Range("E842").Value = 13
$ws.Range("F842").Value = 100112040
$ws.Range("G842").Value = "Cilantro"
$ws.Range("H842").Value = "Sin especificar"
$ws.Range("I842").Value = "Primera"
$ws.Range("J842").Value = 490
$ws.Range("K842").Value = 9000
$ws.Range("L842").Value = 10000
$ws.Range("M842").Value = 9469
$ws.Range("N842").Value = "$/caja 36 atados"
$ws.Range("O842").Value = "Región Metropolitana"
$ws.Range("P842").Value = 263
$ws.Range("Q842").Value = 36
$ws.Range("R842").Value = "Hortaliza"

# New row 843: "Primera" quality, "$/docena de atados" unit, date 2022-07-27
$ws.Range("A843").Value = 6
$ws.Range("B843").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C843").Value = "Metropolitana"
$ws.Range("D843").Value = "2022-07-27"
$ws.Range("E843").Value = 13
$ws.Range("F843").Value = 100112040
$ws.Range("G843").Value = "Cilantro"
$ws.Range("H843").Value = "Sin especificar"
$ws.Range("I843").Value = "Primera"
$ws.Range("J843").Value = 280
$ws.Range("K843").Value = 19000
$ws.Range("L843").Value = 20000
$ws.Range("M843").Value = 19464
$ws.Range("N843").Value = "$/docena de atados"
$ws.Range("O843").Value = "Región Metropolitana"
$ws.Range("P843").Value = 6488
$ws.Range("Q843").Value = 3
$ws.Range("R843").Value = "Hortaliza"
